$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '29.374.15'
Set-TextValue $ws.Range('E2') '  -0.10%  '

Set-TextValue $ws.Range('D3') '1.847.92'
Set-TextValue $ws.Range('E3') '  -0.09%  '

Set-TextValue $ws.Range('D4') '0.9982'
Set-TextValue $ws.Range('E4') '  -0.13%  '

Set-TextValue $ws.Range('D5') '240.22'
Set-TextValue $ws.Range('E5') '  -0.16%  '

Set-TextValue $ws.Range('D6') '0.6270'
Set-TextValue $ws.Range('E6') '  -0.48%  '

Set-TextValue $ws.Range('D7') '0.9996'
Set-TextValue $ws.Range('E7') '  -0.07%  '

Set-TextValue $ws.Range('D8') '0.07616'
Set-TextValue $ws.Range('E8') '  -0.80%  '

Set-TextValue $ws.Range('D9') '0.2902'
Set-TextValue $ws.Range('E9') '  -1.33%  '

Set-TextValue $ws.Range('D10') '24.73'
Set-TextValue $ws.Range('E10') '  +0.92%  '

Set-TextValue $ws.Range('D11') '0.07735'
Set-TextValue $ws.Range('E11') '  -0.17%  '

Set-TextValue $ws.Range('D12') '5.028'
Set-TextValue $ws.Range('E12') '  +0.24%  '

Set-TextValue $ws.Range('D13') '0.6790'
Set-TextValue $ws.Range('E13') '  -0.17%  '

Set-TextValue $ws.Range('D14') '0.00001053'
Set-TextValue $ws.Range('E14') '  -3.58%  '

Set-TextValue $ws.Range('D15') '83.03'
Set-TextValue $ws.Range('E15') '  -0.67%  '

Set-TextValue $ws.Range('D16') '6.155'
Set-TextValue $ws.Range('E16') '  +0.13%  '

Set-TextValue $ws.Range('D17') '29.385.38'
Set-TextValue $ws.Range('E17') '  -0.09%  '

Set-TextValue $ws.Range('D18') '227.61'
Set-TextValue $ws.Range('E18') '  -0.88%  '

Set-TextValue $ws.Range('D19') '12.34'
Set-TextValue $ws.Range('E19') '  -0.93%  '

Set-TextValue $ws.Range('D20') '0.9992'
Set-TextValue $ws.Range('E20') '  -0.11%  '

Set-TextValue $ws.Range('D21') '7.479'
Set-TextValue $ws.Range('E21') '  +0.49%  '

Set-TextValue $ws.Range('D22') '0.9988'
Set-TextValue $ws.Range('E22') '  -0.17%  '

Set-TextValue $ws.Range('D23') '158.52'
Set-TextValue $ws.Range('E23') '  +0.79%  '

Set-TextValue $ws.Range('D24') '0.1385'
Set-TextValue $ws.Range('E24') '  -0.38%  '

Set-TextValue $ws.Range('D25') '8.408'
Set-TextValue $ws.Range('E25') '  +0.53%  '

Set-TextValue $ws.Range('D26') '17.67'
Set-TextValue $ws.Range('E26') '  -0.03%  '

Set-TextValue $ws.Range('D27') '1.408'
Set-TextValue $ws.Range('E27') '  +7.44%  '

Set-TextValue $ws.Range('D28') '1.462'
Set-TextValue $ws.Range('E28') '  -0.39%  '

Set-TextValue $ws.Range('D29') '0.05599'

Set-TextValue $ws.Range('D30') '4.109'
Set-TextValue $ws.Range('E30') '  -0.10%  '

Set-TextValue $ws.Range('E31') '  +0.45%  '

Set-TextValue $ws.Range('B32') 'ARBITRUM'
Set-TextValue $ws.Range('C32') 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range('D32') '1.163'
Set-TextValue $ws.Range('E32') '  +0.39%  '

Set-TextValue $ws.Range('B33') 'LidoDAOToken'
Set-TextValue $ws.Range('C33') 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue $ws.Range('D33') '1.836'
Set-TextValue $ws.Range('E33') '  -0.85%  '

Set-TextValue $ws.Range('D34') '0.6998'
Set-TextValue $ws.Range('E34') '  -1.43%  '

Set-TextValue $ws.Range('E35') '  +0.15%  '

Set-TextValue $ws.Range('B36') 'Maker'
Set-TextValue $ws.Range('C36') 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws.Range('D36') '1.228.85'
Set-TextValue $ws.Range('E36') '  -0.28%  '

Set-TextValue $ws.Range('B37') 'VeChain'
Set-TextValue $ws.Range('C37') 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range('D37') '0.01799'
Set-TextValue $ws.Range('E37') '  +0.06%  '

Set-TextValue $ws.Range('E38') '  -1.82%  '

Set-TextValue $ws.Range('D39') '6.384'
Set-TextValue $ws.Range('E39') '  -1.41%  '

Set-TextValue $ws.Range('D40') '0.9010'
Set-TextValue $ws.Range('E40') '  -1.48%  '

Set-TextValue $ws.Range('D41') '0.9997'
Set-TextValue $ws.Range('E41') '  -0.04%  '

Set-TextValue $ws.Range('D42') '101.36'
Set-TextValue $ws.Range('E42') '  -0.14%  '

Set-TextValue $ws.Range('D43') '65.84'
Set-TextValue $ws.Range('E43') '  -0.57%  '

Set-TextValue $ws.Range('D44') '7.207'
Set-TextValue $ws.Range('E44') '  +0.60%  '

Set-TextValue $ws.Range('B45') 'TheSandbox'
Set-TextValue $ws.Range('C45') 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue $ws.Range('D45') '0.3997'
Set-TextValue $ws.Range('E45') '  -0.44%  '

Set-TextValue $ws.Range('B46') 'EnergySwap'
Set-TextValue $ws.Range('C46') 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range('D46') '9.022'
Set-TextValue $ws.Range('E46') '  -0.02%  '

Set-TextValue $ws.Range('B47') 'RenderToken'
Set-TextValue $ws.Range('C47') 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range('D47') '1.676'
Set-TextValue $ws.Range('E47') '  -0.80%  '

Set-TextValue $ws.Range('B48') 'Algorand'
Set-TextValue $ws.Range('C48') 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Range('D48') '0.1137'
Set-TextValue $ws.Range('E48') '  +1.09%  '

Set-TextValue $ws.Range('B49') 'Cronos'
Set-TextValue $ws.Range('C49') 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Range('D49') '0.05703'
Set-TextValue $ws.Range('E49') '  -0.20%  '

Set-TextValue $ws.Range('B50') 'Mantle'
Set-TextValue $ws.Range('C50') 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue $ws.Range('D50') '0.4623'
Set-TextValue $ws.Range('E50') '  -0.08%  '

Set-TextValue $ws.Range('B51') 'NEARProtocol'
Set-TextValue $ws.Range('C51') 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Range('D51') '1.340'
Set-TextValue $ws.Range('E51') '  -0.80%  '
